$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("classFields")

$ws.Range("B2").Value = "name"
$ws.Range("C2").Value = "private"
$ws.Range("D2").Value = "java.lang.String"
$ws.Range("B4").Value = "ADMIN"
$ws.Range("B5").Value = "BLOGGER"
$ws.Range("C5").Value = "public"
$ws.Range("D5").Value = "org.andante.config.security.role.KeycloakRole"
$ws.Range("B7").Value = "allowedMethods"
$ws.Range("B10").Value = "allowedOrigins"
$ws.Range("B11").Value = "ROLES"
$ws.Range("B12").Value = "REALM_ACCESS"
$ws.Range("B13").Value = "allowedHeaders"
$ws.Range("B14").Value = "exposedHeaders"
$ws.Range("B15").Value = "allowedMethods"
$ws.Range("B16").Value = "disabledSecurityEndpoints"
$ws.Range("D16").Value = "java.util.List"
$ws.Range("B17").Value = "jwkSetUri"
$ws.Range("D17").Value = "java.lang.String"
$ws.Range("B18").Value = "allowedOrigins"
$ws.Range("D18").Value = "java.util.List"
$ws.Range("B19").Value = "keycloakRealmRoleConverter"
$ws.Range("D19").Value = "org.andante.config.security.converter.KeycloakRealmRoleConverter"
